$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto snapshot values.
# NumberFormat "@" (Text) forces Excel to keep these as text strings rather than
# auto-converting plain-looking numbers (e.g. "253.68") into numeric cells,
# matching the original workbook which stores every value as inline text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.322.70"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.26"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.721"
$ws.Range("E5").Value = "  +9.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.68"
$ws.Range("E6").Value = "  +3.63%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").Value = "  -1.98%  "

$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.26"
$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("E11").Value = "  +4.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0989"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.191.06"
$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.62"
$ws.Range("E14").Value = "  +4.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.719"
$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.912.74"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.323.46"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.16"
$ws.Range("E19").Value = "  +2.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0847"
$ws.Range("E20").Value = "  +2.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.84"
$ws.Range("E21").Value = "  +1.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.00"
$ws.Range("E22").Value = "  +3.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.07"
$ws.Range("E23").Value = "  +3.84%  "

$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +4.19%  "

$ws.Range("E26").Value = "  +3.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.91"
$ws.Range("E27").Value = "  -2.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.61"
$ws.Range("E28").Value = "  +1.73%  "

$ws.Range("E29").Value = "  +1.02%  "

$ws.Range("E30").Value = "  +3.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.126.49"
$ws.Range("E31").Value = "  +19.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.34"
$ws.Range("E32").Value = "  +4.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.99"
$ws.Range("E33").Value = "  +14.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.62"
$ws.Range("E34").Value = "  +21.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0583"
$ws.Range("E35").Value = "  +2.63%  "

$ws.Range("E36").Value = "  +1.70%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("E38").Value = "  -2.93%  "

$ws.Range("E39").Value = "  -1.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.33"
$ws.Range("E40").Value = "  +6.61%  "

$ws.Range("E41").Value = "  +3.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.11"
$ws.Range("E42").Value = "  +6.94%  "

$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0640"
$ws.Range("E44").Value = "  -1.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.337.54"
$ws.Range("E45").Value = "  -0.55%  "

$ws.Range("E46").Value = "  +1.72%  "

$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.78"
$ws.Range("E48").Value = "  +2.90%  "

$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.30"
$ws.Range("E50").Value = "  +25.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.08"
$ws.Range("E51").Value = "  -7.63%  "
